$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''33.570.16'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  -1.18%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''1.760.22'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -1.25%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  +0.42%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''223.24'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +0.92%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''0.540'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  -2.13%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').Value = '''  +0.37%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = '''31.66'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '''  +1.07%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = '''0.285'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '''  -0.49%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''0.0682'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  -3.80%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('E11').Value = '''  +1.73%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = '''2.016.76'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  -1.08%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = '''11.10'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '''  +6.12%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = '''1.784.26'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '''  +0.20%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''33.629.06'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  -0.97%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''0.605'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  -3.32%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''4.08'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  -3.12%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''66.29'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  -2.49%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''0.0₃0765'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '''  -1.74%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''235.78'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  -3.62%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('E21').Value = '''  +0.34%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''10.49'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  -1.65%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''4.00'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''  -1.71%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''2.05'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  -1.89%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''158.74'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +0.83%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''16.01'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''  -2.23%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = '''6.97'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '''  -0.29%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = '''0.111'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''  -0.79%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('E29').Value = '''  +0.54%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('E30').Value = '''  +1.50%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('E31').Value = '''  -2.61%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''3.57'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  -3.31%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').Value = '''  -0.38%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''1.76'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  -2.36%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = '''1.377.29'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '''  -1.61%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = '''0.648'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '''  +1.61%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = '''1.02'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '''  -2.48%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('E38').Value = '''  -1.58%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D39').Value = '''2.36'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '''  +0.96%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''2.20'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  +4.97%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = '''0.901'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''  -3.43%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = '''77.20'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  -2.55%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('E43').Value = '''  -2.70%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = '''13.34'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  +13.76%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = '''0.0₆0139'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  +16.27%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('E46').Value = '''  +4.44%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = '''0.0498'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '''  +1.52%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = '''106.76'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '''  +1.75%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = '''5.79'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '''  -2.41%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = '''1.917.40'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '''  -0.66%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('E51').Value = '''  +0.51%  '
$ws.Range('E51').Style = "Normal"
